# Weekly update of the Achicoria (Vega Modelo de Temuco) price series:
# insert two new rows of data into the existing table, pushing the
# older rows down, and fill in the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the newest observation right after the header/earlier rows,
#     at row 54 (everything from old row 54 down shifts to row 55+). ---
$ws.Rows.Item(54).Insert()

$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 45007
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = 100112010
$ws.Cells.Item(54, 7).Value = "Achicoria"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 25
$ws.Cells.Item(54, 11).Value = 10000
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 10000
$ws.Cells.Item(54, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value = 556
$ws.Cells.Item(54, 17).Value = 18
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# --- Insert a second new observation at (the now shifted) row 70, just
#     before the former last row (which moves down to row 71). ---
$ws.Rows.Item(70).Insert()

$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 45008
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = 100112010
$ws.Cells.Item(70, 7).Value = "Achicoria"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 65
$ws.Cells.Item(70, 11).Value = 10000
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = 10000
$ws.Cells.Item(70, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(70, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 16).Value = 556
$ws.Cells.Item(70, 17).Value = 18
$ws.Cells.Item(70, 18).Value = "Hortaliza"
